$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/rank refresh, matching the GitHub Actions data pull
# on 2022-12-24 17:45:58 UTC. Each entry is the new literal cell value for
# the "cryptos" sheet (Coin / Link / Price / Volume(1h) columns).
$edits = @(
    @{ Cell = "D2"; Value = '244.60' },
    @{ Cell = "D3"; Value = '21.83' },
    @{ Cell = "D4"; Value = '5.387' },
    @{ Cell = "D7"; Value = '0.8152' },
    @{ Cell = "D8"; Value = '0.9309' },
    @{ Cell = "B9"; Value = 'One' },
    @{ Cell = "C9"; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' },
    @{ Cell = "D9"; Value = '0.0005942' },
    @{ Cell = "E9"; Value = '8OneONE' },
    @{ Cell = "B10"; Value = 'WazirX' },
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = "D10"; Value = '0.1432' },
    @{ Cell = "E10"; Value = '9WazirXWRX' },
    @{ Cell = "B11"; Value = 'MandalaExchangeToken' },
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = "D11"; Value = '0.07361' },
    @{ Cell = "E11"; Value = '10MandalaExchangeTokenMDX' },
    @{ Cell = "B12"; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = "D12"; Value = '0.03452' },
    @{ Cell = "E12"; Value = '11LiechtensteinCryptoassetsExchangeLCX' },
    @{ Cell = "B13"; Value = 'BitrueCoin' },
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = "D13"; Value = '0.03061' },
    @{ Cell = "E13"; Value = '12BitrueCoinBTR' },
    @{ Cell = "B14"; Value = 'BitMartToken' },
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = "D14"; Value = '0.09407' },
    @{ Cell = "E14"; Value = '13BitMartTokenBMX' },
    @{ Cell = "B15"; Value = 'MCDex' },
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' },
    @{ Cell = "D15"; Value = '4.001' },
    @{ Cell = "E15"; Value = '14MCDexMCB' },
    @{ Cell = "B16"; Value = 'BitForexToken' },
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = "D16"; Value = '0.001594' },
    @{ Cell = "E16"; Value = '15BitForexTokenBF' },
    @{ Cell = "B17"; Value = 'CoinExToken' },
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' },
    @{ Cell = "D17"; Value = '0.04808' },
    @{ Cell = "E17"; Value = '16CoinExTokenCET' },
    @{ Cell = "D18"; Value = '0.005563' },
    @{ Cell = "D19"; Value = '0.004154' },
    @{ Cell = "D20"; Value = '0.0009867' },
    @{ Cell = "D40"; Value = '0.04013' },
    @{ Cell = "D41"; Value = '0.006412' },
    @{ Cell = "E41"; Value = '40KickTokenKICKBestin24h' },
    @{ Cell = "D42"; Value = '0.1074' },
    @{ Cell = "D44"; Value = '0.006668' },
    @{ Cell = "E44"; Value = '43LocalTradersLCT' },
    @{ Cell = "D45"; Value = '0.00005239' },
    @{ Cell = "D48"; Value = '0.002552' }
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    $value = $edit.Value

    # Price-like values ("244.60", "0.006412", ...) are stored as literal
    # text in this sheet (the same way the source CSV -> XLSX pipeline
    # writes them), not as numbers. A leading apostrophe forces Excel to
    # keep the exact text instead of re-parsing it into a Double (which
    # would both change the stored type and introduce float rounding such
    # as 244.6 for "244.60"). Resetting the style back to "Normal"
    # afterwards clears the quote-prefix flag Excel sets on the cell so no
    # extra formatting is introduced.
    $isNumericLooking = $value -match '^-?\d+(\.\d+)?$'
    if ($isNumericLooking) {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
    $cell.Style = "Normal"
}
